# Applies the "Add files via upload" edit to Trivago.xlsx:
#  - Adds a prediction/label column (B) to the "Teste" sheet, with lowercase
#    "irrelevante"/"relevante" labels for each of the 200 tweet rows (2-201).
#  - Sets column A width on "Teste" to match the author's saved view.
#  - Updates the active sheet / selection state so that "Teste" becomes the
#    selected tab, matching the workbook state captured in the commit.

$wb = $excel.ActiveWorkbook
$wsTrain = $wb.Worksheets.Item("Treinamento")
$wsTest  = $wb.Worksheets.Item("Teste")

# Compact encoding of the label for rows 2..201 of the "Teste" sheet:
# 'i' => "irrelevante", 'r' => "relevante"
$labelCodes = "iiiirrriirriiiiirririiiiriiiiiriiriiriiiiiriiiiiiiiiiiirrririiiirriirirriiirrriiiirriiririiiiiiiririiiiiiiiiiiriiiiiiiirriiiiiiiiriiririiririiiiiiirriiriiirrriiiiiiiiiiiiiiriiiiririiiriiiiiiiiiiiririi"

for ($i = 0; $i -lt $labelCodes.Length; $i++) {
    $row = $i + 2
    $code = $labelCodes.Substring($i, 1)
    if ($code -eq "i") {
        $label = "irrelevante"
    } else {
        $label = "relevante"
    }
    $wsTest.Cells.Item($row, 2).Value = $label
}

# Column A width on the "Teste" sheet (author widened it to fit tweet text).
$wsTest.Columns.Item(1).ColumnWidth = 118.75

# Recreate the saved window/selection state: "Treinamento" scrolled to A168
# with B179 selected, "Teste" scrolled to A181 with B202 selected, and
# "Teste" left as the active (selected) sheet/tab.
$wsTrain.Activate()
$excel.ActiveWindow.ScrollRow = 168
$excel.ActiveWindow.ScrollColumn = 1
$wsTrain.Range("B179").Select() | Out-Null

$wsTest.Activate()
$excel.ActiveWindow.ScrollRow = 181
$excel.ActiveWindow.ScrollColumn = 1
$wsTest.Range("B202").Select() | Out-Null
